# Edit script implementing the commit:
# "Elimina antiguos EC y agrega nuevos y modifica Antigua BD"
# (Delete old account statements and add new ones, and modify the old BD/period)
#
# Summary of changes applied to Hoja1:
#  - Insert a new row for worker 72138037 / FRANKLIN BENJAMIN ESCAMILLA CONTRERAS
#    for period 2509 (right after his existing 2508 row), keeping the same
#    Valor Mora / Salario Basico values.
#  - Update the "Periodo Mora" (column E) of every other worker row from
#    period 2508 to the new period 2509.
#  - Update the total "VALOR MORA" cell and "Cant. Periodos" counter.
#  - Center-align the "Periodo Mora" column in the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. Update summary header fields
# ---------------------------------------------------------------------
# VALOR MORA total
$ws.Range("E11").Value = 1278089
# Cant. Periodos
$ws.Range("F13").Value = 2

# ---------------------------------------------------------------------
# 2. Insert a new data row right below row 16 (the first worker),
#    duplicating that worker for the new period 2509.
# ---------------------------------------------------------------------
$ws.Rows.Item(17).Insert()

# Copy formatting from the row above (16) into the newly inserted row (17)
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "72138037"
$ws.Range("D17").Value = "FRANKLIN BENJAMIN ESCAMILLA CONTRERAS"
$ws.Range("E17").Value = "2509"
$ws.Range("F17").Value = 481354
$ws.Range("G17").Value = 12033840

# ---------------------------------------------------------------------
# 3. Update period (column E) for the remaining workers: 2508 -> 2509
#    (rows shifted down by 1 because of the inserted row above)
# ---------------------------------------------------------------------
$ws.Range("E18").Value = "2509"
$ws.Range("E19").Value = "2509"
$ws.Range("E20").Value = "2509"
$ws.Range("E21").Value = "2509"

# ---------------------------------------------------------------------
# 4. Center the "Periodo Mora" column across the whole data table
# ---------------------------------------------------------------------
$ws.Range("E16:E21").HorizontalAlignment = -4108
